$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look like plain numbers (e.g. "512.44").
# Force them to remain text by pre-formatting each as Text ("@") before assigning,
# otherwise Excel would silently coerce them into numeric values. (Looping per-cell
# because applying NumberFormat to a multi-area union range only affects the first area.)
$textCells = @("D5", "D6", "D7", "D9", "D10", "D14", "D18", "D19", "D20", "D21", "D22", "D24", "D26", "D28", "D29", "D32", "D33", "D34", "D38", "D40", "D46", "D47", "D50")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '56.829.37'
$ws.Range("E2").Value = '  +2.54%  '
$ws.Range("D3").Value = '3.000.83'
$ws.Range("E3").Value = '  +2.10%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '512.44'
$ws.Range("E5").Value = '  +5.21%  '
$ws.Range("D6").Value = '138.57'
$ws.Range("E6").Value = '  +5.73%  '
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.21%  '
$ws.Range("E8").Value = '  +3.95%  '
$ws.Range("D9").Value = '7.52'
$ws.Range("E9").Value = '  +5.49%  '
$ws.Range("D10").Value = '0.108'
$ws.Range("E10").Value = '  +7.68%  '
$ws.Range("E11").Value = '  +3.23%  '
$ws.Range("E12").Value = '  +2.57%  '
$ws.Range("D13").Value = '3.514.00'
$ws.Range("E13").Value = '  +1.91%  '
$ws.Range("D14").Value = '25.65'
$ws.Range("E14").Value = '  +4.90%  '
$ws.Range("E15").Value = '  +12.46%  '
$ws.Range("D16").Value = '56.833.99'
$ws.Range("E16").Value = '  +2.45%  '
$ws.Range("D17").Value = '3.001.99'
$ws.Range("E17").Value = '  +1.85%  '
$ws.Range("D18").Value = '5.90'
$ws.Range("E18").Value = '  +6.18%  '
$ws.Range("D19").Value = '12.53'
$ws.Range("E19").Value = '  +3.66%  '
$ws.Range("D20").Value = '7.84'
$ws.Range("E20").Value = '  +5.09%  '
$ws.Range("D21").Value = '326.66'
$ws.Range("E21").Value = '  +3.53%  '
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("E23").Value = '  +5.76%  '
$ws.Range("D24").Value = '63.26'
$ws.Range("E24").Value = '  +5.77%  '
$ws.Range("E25").Value = '  +6.97%  '
$ws.Range("D26").Value = '0.995'
$ws.Range("E26").Value = '  -0.61%  '
$ws.Range("D27").Value = '0.0₃0914'
$ws.Range("E27").Value = '  +9.97%  '
$ws.Range("D28").Value = '6.62'
$ws.Range("E28").Value = '  +3.04%  '
$ws.Range("D29").Value = '7.07'
$ws.Range("E29").Value = '  +8.66%  '
$ws.Range("E30").Value = '  +6.91%  '
$ws.Range("E31").Value = '  +6.91%  '
$ws.Range("D32").Value = '20.55'
$ws.Range("E32").Value = '  +6.32%  '
$ws.Range("D33").Value = '157.19'
$ws.Range("E33").Value = '  +6.38%  '
$ws.Range("D34").Value = '4.57'
$ws.Range("E34").Value = '  +5.27%  '
$ws.Range("E35").Value = '  +0.71%  '
$ws.Range("E36").Value = '  -1.22%  '
$ws.Range("E37").Value = '  +4.51%  '
$ws.Range("D38").Value = '23.89'
$ws.Range("E38").Value = '  +4.11%  '
$ws.Range("D39").Value = '3.032.45'
$ws.Range("E39").Value = '  +2.02%  '
$ws.Range("D40").Value = '37.09'
$ws.Range("E40").Value = '  +2.74%  '
$ws.Range("E41").Value = '  -0.12%  '
$ws.Range("D42").Value = '2.286.00'
$ws.Range("E42").Value = '  +8.50%  '
$ws.Range("E43").Value = '  +3.47%  '
$ws.Range("E44").Value = '  +4.59%  '
$ws.Range("E45").Value = '  +3.43%  '
$ws.Range("D46").Value = '0.999'
$ws.Range("E46").Value = '  +0.73%  '
$ws.Range("D47").Value = '1.95'
$ws.Range("E47").Value = '  +10.76%  '
$ws.Range("E48").Value = '  +3.00%  '
$ws.Range("E49").Value = '  +6.39%  '
$ws.Range("D50").Value = '19.19'
$ws.Range("E50").Value = '  +0.57%  '
$ws.Range("E51").Value = '  +4.96%  '
